$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B6 previously held the text "01/06/2021  Taller". Replace it with the
# actual date value (2021-06-01) formatted as a date, keeping the existing
# border/alignment formatting.
$ws.Range("B6").Value = Get-Date -Year 2021 -Month 6 -Day 1 -Hour 0 -Minute 0 -Second 0

# Move the active selection from D6 to B7.
$ws.Range("B7").Select() | Out-Null
